$d = $word.ActiveDocument

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $fullXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>$innerXml</w:p>"
    $r.InsertXML($fullXml)
}

# Paragraph 5 (1-indexed) / index 4: "As you can see..."
Set-ParagraphXml 5 '<w:r><w:t xml:space="preserve">As you can see, the more mass an object </w:t></w:r><w:r><w:t>has</w:t></w:r><w:r><w:t>, the less it will accelerate from the net force.</w:t></w:r>'

# Paragraph 9 (1-indexed) / index 8: vector sentence
Set-ParagraphXml 9 '<w:r><w:t>A vector is composed of a value for each axis. In our case, the x and y values. These values can also tell us the direction, and the magnitude (for example: speed).</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r>'

# Paragraph 12 (1-indexed) / index 11: wheel intro paragraph - merge last two runs, drop bookmark (re-added later)
Set-ParagraphXml 12 '<w:r><w:t xml:space="preserve">A </w:t></w:r><w:r><w:t>b</w:t></w:r><w:r><w:t xml:space="preserve">rave </w:t></w:r><w:r><w:t xml:space="preserve">sir </w:t></w:r><w:r><w:t xml:space="preserve">knight </w:t></w:r><w:r><w:t>has</w:t></w:r><w:r><w:t xml:space="preserve"> brought us a wheel to demonstrate. </w:t></w:r><w:r><w:t>The wheel is enchanted with magical trails to allow us to observe its motion.</w:t></w:r>'

# Paragraph 19 (1-indexed) / index 18: "Around here..." + bookmark to be added
Set-ParagraphXml 19 '<w:r><w:t>Around here, the distance between each trail is the same. This tells us that the net force on the wheel equals zero.</w:t></w:r>'

# Paragraph 22 (1-indexed) / index 21: "At this point the trails..."
Set-ParagraphXml 22 '<w:r><w:t>At this point the trails are going down. We can observe that the only force acting on the wheel is the gravity.</w:t></w:r>'

# Move the _GoBack bookmark: collapsed, right after "Around here" in paragraph 19 (1-indexed)
$p19 = $d.Paragraphs.Item(19)
$bmPos = $p19.Range.Start + ("Around here").Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Move lastRenderedPageBreak: remove from "(enable graph button...)" paragraph, add to preceding "(hide highlight)" paragraph
# First: rewrite the "(enable graph button...)" paragraph so the lastRenderedPageBreak + proofErr go away, and "Show" -> "show"
$d.Content.Find.Execute("(enable graph button, Show pop-up on it:)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(enable graph button, show pop-up on it:)", 2)

# Then: add lastRenderedPageBreak to the "(hide highlight)" paragraph right before it (index 22, 1-indexed 23)
Set-ParagraphXml 23 '<w:r><w:lastRenderedPageBreak/><w:t>(hide highlight)</w:t></w:r>'

Write-Output $d.Paragraphs.Count
for ($i=1; $i -le 25; $i++) {
    Write-Output "$i : [$($d.Paragraphs.Item($i).Range.Text)]"
}
